$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# The perf-test results sheet currently has two rows of "WINDOWS" samples
# (row 6 = RequestResponse, row 7 = DBExecuting). This change adds the
# matching "LINUX" samples for the same run, inserted *above* the existing
# WINDOWS rows so the WINDOWS/LINUX pairing pattern used in rows 2-5 is
# restored (avg - N * stdevp maybe < 0 fix surfaced while re-running stats).
# ---------------------------------------------------------------------------

# Insert a blank row at 6 (existing row 6 -> 7) and another blank row at 8
# (existing row 7, now at 8, -> 9).
$ws.Rows.Item(6).Insert()
$ws.Rows.Item(8).Insert()

# Copy the per-cell formatting (fill/border/number format) from the rows
# that were pushed down into the freshly inserted blank rows, so the new
# rows look like the rest of the shaded table instead of plain default
# cells.
$ws.Range("A7:S7").Copy()
$ws.Range("A6:S6").PasteSpecial(-4122)
$ws.Range("A9:S9").Copy()
$ws.Range("A8:S8").PasteSpecial(-4122)
$ws.Range("A1").Select()

# -- New row 6: RequestResponse / LINUX ------------------------------------
$ws.Range("A6").Value2 = "RequestResponse"
$ws.Range("B6").Value2 = "2020-02-03 02:26:09.000"
$ws.Range("C6").Value2 = "LINUX"
$ws.Range("D6").Value2 = "Unix 5.3.0.26"
$ws.Range("E6").Value2 = 210
$ws.Range("F6").Value2 = 11.3048
$ws.Range("G6").Value2 = 27.370142999999999
$ws.Range("H6").Value2 = 788.68650000000002
$ws.Range("I6").Value2 = 14.844799999999999
$ws.Range("J6").Value2 = 17.833967568693101
$ws.Range("K6").Value2 = 10000
$ws.Range("L6").Value2 = 9776
$ws.Range("M6").Value2 = 92.92
$ws.Range("N6").Value2 = 0
$ws.Range("N6").Font.Color = 255
$ws.Range("O6").Value2 = 63.038078137386201
$ws.Range("P6").Value2 = 10001
$ws.Range("Q6").Value2 = 20000
$ws.Range("R6").Value2 = "2020-02-03 02:26:31.000"
$ws.Range("S6").Value2 = "2020-02-03 02:30:01.000"

# -- New row 8: DBExecuting / LINUX -----------------------------------------
$ws.Range("A8").Value2 = "DBExecuting"
$ws.Range("B8").Value2 = "2020-02-03 02:26:09.000"
$ws.Range("C8").Value2 = "LINUX"
$ws.Range("D8").Value2 = "Unix 5.3.0.26"
$ws.Range("E8").Value2 = 210
$ws.Range("F8").Value2 = 7.8936999999999999
$ws.Range("G8").Value2 = 16.373159999999999
$ws.Range("H8").Value2 = 450.78440000000001
$ws.Range("I8").Value2 = 11.305899999999999
$ws.Range("J8").Value2 = 9.55310322735399
$ws.Range("K8").Value2 = 10000
$ws.Range("L8").Value2 = 9752
$ws.Range("M8").Value2 = 91.38
$ws.Range("N8").Value2 = 0
$ws.Range("N8").Font.Color = 255
$ws.Range("O8").Value2 = 35.479366454708
$ws.Range("P8").Value2 = 10001
$ws.Range("Q8").Value2 = 20000
$ws.Range("R8").Value2 = "2020-02-03 02:26:31.000"
$ws.Range("S8").Value2 = "2020-02-03 02:30:01.000"

# Restore the selection that Excel recorded the sheet with after the edit.
$ws.Range("M32").Select()

Write-Output "done"
